# Applies the profit-table value updates captured in the commit diff.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 395
$ws.Range("I31").Value = 395
$ws.Range("K31").Value = 1185
$ws.Range("M31").Value = -955

$ws.Range("H34").Value = 11735.154
$ws.Range("I34").Value = 11835.917
$ws.Range("J34").Value = 10526
$ws.Range("K34").Value = 11835.917
$ws.Range("L34").Value = 10526
$ws.Range("M34").Value = -11632.917
$ws.Range("N34").Value = -10932

$ws.Range("H36").Value = 11735.154
$ws.Range("I36").Value = 11835.917
$ws.Range("J36").Value = 10526
$ws.Range("K36").Value = 11835.917
$ws.Range("L36").Value = 10526
$ws.Range("M36").Value = -11120.917
$ws.Range("N36").Value = -11956

$ws.Range("H38").Value = 1204.5454
$ws.Range("I38").Value = 1204.5454
$ws.Range("K38").Value = 3613.6362
$ws.Range("M38").Value = -3241.6362

$ws.Range("H64").Value = 7812.25
$ws.Range("J64").Value = 8215.143
$ws.Range("L64").Value = 8215.143
$ws.Range("N64").Value = -8711.143

$ws.Range("H67").Value = 7812.25
$ws.Range("J67").Value = 8215.143
$ws.Range("L67").Value = 8215.143
$ws.Range("N67").Value = -9931.143

$ws.Range("H92").Value = 683.1818
$ws.Range("I92").Value = 112.77778
$ws.Range("K92").Value = 112.77778
$ws.Range("M92").Value = 1135.22222

$ws.Range("H100").Value = 8809.066000000001
$ws.Range("J100").Value = 13444.556
$ws.Range("L100").Value = 13444.556
$ws.Range("N100").Value = -14526.556

$ws.Range("H131").Value = 3022.8948
$ws.Range("I131").Value = 2158.3845
$ws.Range("J131").Value = 4896
$ws.Range("K131").Value = 6475.1535
$ws.Range("L131").Value = 14688
$ws.Range("M131").Value = -1435.1535
$ws.Range("N131").Value = -24768

$ws.Range("H141").Value = 2373.7144
$ws.Range("I141").Value = 2436.5557
$ws.Range("K141").Value = 7309.6671
$ws.Range("M141").Value = -2129.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5257.52
$ws.Range("I45").Value = 1996
$ws.Range("K45").Value = 1996
$ws.Range("M45").Value = -1619

$ws.Range("H88").Value = 2317.3333
$ws.Range("J88").Value = 2476
$ws.Range("L88").Value = 2476
$ws.Range("N88").Value = -3288

$ws.Range("H91").Value = 2317.3333
$ws.Range("J91").Value = 2476
$ws.Range("L91").Value = 2476
$ws.Range("N91").Value = -5284

$ws.Range("H124").Value = 53747
$ws.Range("J124").Value = 53747
$ws.Range("L124").Value = 53747
$ws.Range("N124").Value = -63567

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

$ws.Range("H20").Value = 2524.2
$ws.Range("I20").Value = 1934.1
$ws.Range("K20").Value = 1934.1
$ws.Range("M20").Value = -1687.1

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H125").Value = 74995
$ws.Range("J125").Value = 74995
$ws.Range("L125").Value = 74995
$ws.Range("N125").Value = -84835

$ws.Range("H134").Value = 4339.8
$ws.Range("I134").Value = 2841.4167
$ws.Range("K134").Value = 8524.250100000001
$ws.Range("M134").Value = -5989.250100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 3000
$ws.Range("I37").Value = 3000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 3000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2893
$ws.Range("N37").ClearContents()

$ws.Range("H58").Value = 181867.86
$ws.Range("I58").Value = 296293
$ws.Range("J58").Value = 5029
$ws.Range("K58").Value = 296293
$ws.Range("L58").Value = 5029
$ws.Range("M58").Value = -296090
$ws.Range("N58").Value = -5435

$ws.Range("H105").Value = 1285.4348
$ws.Range("J105").Value = 1286.1428
$ws.Range("L105").Value = 1286.1428
$ws.Range("N105").Value = -4780.1428

$ws.Range("H132").Value = 3620.6956
$ws.Range("I132").Value = 2540.6667
$ws.Range("K132").Value = 7622.000100000001
$ws.Range("M132").Value = -5092.000100000001

$ws.Range("H136").Value = 181867.86
$ws.Range("I136").Value = 296293
$ws.Range("J136").Value = 5029
$ws.Range("K136").Value = 888879
$ws.Range("L136").Value = 15087
$ws.Range("M136").Value = -886329
$ws.Range("N136").Value = -20187

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 322.25
$ws.Range("I103").Value = 233.5
$ws.Range("K103").Value = 700.5
$ws.Range("M103").Value = 178.5

$ws.Range("H121").Value = 834247.5
$ws.Range("I121").Value = 394
$ws.Range("J121").Value = 1429857.1
$ws.Range("K121").Value = 1182
$ws.Range("L121").Value = 4289571.300000001
$ws.Range("M121").Value = 128
$ws.Range("N121").Value = -4292191.300000001

$ws.Range("H128").Value = 207821
$ws.Range("I128").Value = 207821
$ws.Range("K128").Value = 623463
$ws.Range("M128").Value = -618483

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3632.5757
$ws.Range("I97").Value = 4306.963
$ws.Range("J97").Value = 597.8333
$ws.Range("K97").Value = 4306.963
$ws.Range("L97").Value = 597.8333
$ws.Range("M97").Value = -3810.963
$ws.Range("N97").Value = -1589.8333

$ws.Range("H99").Value = 2883
$ws.Range("I99").Value = 2883
$ws.Range("K99").Value = 2883
$ws.Range("M99").Value = -637

$ws.Range("H113").Value = 630328.0600000001
$ws.Range("I113").Value = 671016.9399999999
$ws.Range("J113").Value = 19995
$ws.Range("K113").Value = 671016.9399999999
$ws.Range("L113").Value = 19995
$ws.Range("M113").Value = -668846.9399999999
$ws.Range("N113").Value = -24335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 40689.668
$ws.Range("I13").Value = 30000
$ws.Range("J13").Value = 46034.5
$ws.Range("K13").Value = 30000
$ws.Range("L13").Value = 46034.5
$ws.Range("M13").Value = -29860
$ws.Range("N13").Value = -46314.5

$ws.Range("H14").Value = 2005
$ws.Range("J14").Value = 2005
$ws.Range("L14").Value = 2005
$ws.Range("N14").Value = -2349

$ws.Range("H46").Value = 5381.8823
$ws.Range("I46").Value = 4277.778
$ws.Range("K46").Value = 4277.778
$ws.Range("M46").Value = -4089.778

$ws.Range("H100").Value = 144857.28
$ws.Range("I100").Value = 168500.17
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 168500.17
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -167959.17
$ws.Range("N100").Value = -4082

$ws.Range("H132").Value = 4444.7896
$ws.Range("I132").Value = 2881.5
$ws.Range("J132").Value = 8822
$ws.Range("K132").Value = 8644.5
$ws.Range("L132").Value = 26466
$ws.Range("M132").Value = -6114.5
$ws.Range("N132").Value = -31526

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H8").Value = 18000
$ws.Range("I8").Value = 18000
$ws.Range("K8").Value = 18000
$ws.Range("M8").Value = -17860

$ws.Range("H41").Value = 8597.25
$ws.Range("I41").Value = 5857.6665
$ws.Range("J41").Value = 10241
$ws.Range("K41").Value = 5857.6665
$ws.Range("L41").Value = 10241
$ws.Range("M41").Value = -5467.6665
$ws.Range("N41").Value = -11021

$ws.Range("H46").Value = 93500
$ws.Range("J46").Value = 93500
$ws.Range("L46").Value = 93500
$ws.Range("N46").Value = -93962

$ws.Range("H81").Value = 31798.75
$ws.Range("I81").Value = 55486.75
$ws.Range("K81").Value = 110973.5
$ws.Range("M81").Value = -109912.5

$ws.Range("H84").Value = 31798.75
$ws.Range("I84").Value = 55486.75
$ws.Range("K84").Value = 554867.5
$ws.Range("M84").Value = -549563.5

$ws.Range("H132").Value = 4324.0527
$ws.Range("J132").Value = 6625
$ws.Range("L132").Value = 19875
$ws.Range("N132").Value = -24935

$ws.Range("H134").Value = 93500
$ws.Range("J134").Value = 93500
$ws.Range("L134").Value = 280500
$ws.Range("N134").Value = -285570
